$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the four new worksheets at the end of the workbook, in order:
# CreateWork, Works, WorkinGroup, SubWorks
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCreateWork = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsCreateWork.Name = "CreateWork"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWorks = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsWorks.Name = "Works"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWorkinGroup = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsWorkinGroup.Name = "WorkinGroup"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSubWorks = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsSubWorks.Name = "SubWorks"

# ---------------------------------------------------------------------------
# CreateWork sheet data
# ---------------------------------------------------------------------------
$wsCreateWork.Range("A1").Value = "TCID"
$wsCreateWork.Range("B1").Value = "title"
$wsCreateWork.Range("C1").Value = "description"

$wsCreateWork.Range("A2").Value = 1
$wsCreateWork.Range("B2").Value = "fWorks: Prepare for Testing"
$wsCreateWork.Range("C2").Value = "Đây là dự án cực kỳ quan trọng"

$wsCreateWork.Range("A3").Value = 2
$wsCreateWork.Range("C3").Value = "Đây là dự án cực kỳ quan trọng"

$wsCreateWork.Range("A4").Value = 3
$wsCreateWork.Range("B4").Value = "fWorks: Prepare for Testing"
$wsCreateWork.Range("C4").Value = "Đây là dự án cực kỳ quan trọng"

$wsCreateWork.Columns.Item(1).AutoFit()
$wsCreateWork.Columns.Item(2).AutoFit()
$wsCreateWork.Columns.Item(3).AutoFit()

$wsCreateWork.Range("D10").Select()

# ---------------------------------------------------------------------------
# Works sheet data
# ---------------------------------------------------------------------------
$wsWorks.Range("A1").Value = "TCID"
$wsWorks.Range("B1").Value = "title"

$wsWorks.Range("A2").Value = 1

$wsWorks.Range("A3").Value = 2
$wsWorks.Range("B3").Value = "Thực thi automation đạt chuẩn"

$wsWorks.Range("Q12:Q13").Select()

# ---------------------------------------------------------------------------
# WorkinGroup sheet data
# ---------------------------------------------------------------------------
$wsWorkinGroup.Range("A1").Value = "TCID"
$wsWorkinGroup.Range("B1").Value = "title"

$wsWorkinGroup.Range("A2").Value = 1
$wsWorkinGroup.Range("B2").Value = "Lập plan để thực thi"

$wsWorkinGroup.Range("A3").Value = 2
$wsWorkinGroup.Range("B3").Value = "Lập plan để thực thi"

$wsWorkinGroup.Range("A4").Value = 3

$wsWorkinGroup.Range("A5").Value = 4
$wsWorkinGroup.Range("B5").Value = "Lập plan để thực thi"

$wsWorkinGroup.Range("A6").Value = 5
$wsWorkinGroup.Range("B6").Value = "Lập plan để thực thi"

$wsWorkinGroup.Columns.Item(2).AutoFit()

$wsWorkinGroup.Range("B6").Select()

# ---------------------------------------------------------------------------
# SubWorks sheet data
# (description values are written before title values so the shared-string
# table allocates "Đây là công việc phụ vô cùng quan trọng" before
# "Công việc phụ số 1"/"Công việc phụ số 2", matching the authored order.)
# ---------------------------------------------------------------------------
$wsSubWorks.Range("A1").Value = "TCID"
$wsSubWorks.Range("B1").Value = "title"
$wsSubWorks.Range("C1").Value = "description"

$wsSubWorks.Range("A2").Value = 1
$wsSubWorks.Range("C2").Value = "Đây là công việc phụ vô cùng quan trọng"
$wsSubWorks.Range("B2").Value = "Công việc phụ số 1"

$wsSubWorks.Range("A3").Value = 2
$wsSubWorks.Range("C3").Value = "Đây là công việc phụ vô cùng quan trọng"
$wsSubWorks.Range("B3").Value = "Công việc phụ số 2"

$wsSubWorks.Range("A4").Value = 3
$wsSubWorks.Range("C4").Value = "Đây là công việc phụ vô cùng quan trọng"

$wsSubWorks.Range("A5").Value = 4
$wsSubWorks.Range("C5").Value = "Đây là công việc phụ vô cùng quan trọng"
$wsSubWorks.Range("B5").Value = "Công việc phụ số 1"

$wsSubWorks.Columns.Item(3).AutoFit()

$wsSubWorks.Range("H6").Select()

# SubWorks is the sheet left active/selected, matching tabSelected="1"
$wsSubWorks.Activate()
